$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds numeric-looking strings (e.g. '0.250', '4.30', '26.521.62')
# that must stay literal text exactly as shown (trailing zeros, multi-dot
# thousand groupings, subscript digits, etc.), matching the source feed's
# formatting. Excel's COM layer auto-coerces plain numeric-looking strings
# to real numbers on assignment, which silently drops e.g. a trailing zero
# ('0.250' -> 0.25). Force text via NumberFormat '@' before assigning, then
# ClearFormats() right after so the cell's style index is left untouched
# (matching the unchanged s="" attribute on these cells in the source file) -
# the value itself stays text once committed; only the quote-prefix styling
# flag is cleared.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.521.62'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.627.75'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.487'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.250'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0618'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.97'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0828'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.852.72'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.632.99'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.490.93'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.84'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.72'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.87'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.95'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  -2.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.29'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.61'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0520'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +6.21%  '
$ws.Range("E31").Value = '  +0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.94'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.160.27'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("E37").Value = '  +1.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.805'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.43%  '
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.32'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.21%  '
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.42'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.783'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.764.76'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.18'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("E46").Value = '  +1.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0104'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +6.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.10'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0509'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.409'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.45%  '
$ws.Range("E51").Value = '  +0.14%  '
